# Finish integration of basic GPS code, begin code for movement
#
# The prior commit added a brand new weekly-status sheet ("10-28-13") by
# duplicating the previous week's sheet ("10-21-13"), updating its header
# date, replacing the two task rows with the new GPS-related tasks, and
# moving the "active" tab / selection over to the new sheet.

$wb = $excel.ActiveWorkbook

# The most recent existing weekly sheet is "10-21-13" (last tab).
$prevSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Duplicate it to create next week's sheet, placed right after it.
$prevSheet.Copy([System.Reflection.Missing]::Value, $prevSheet) | Out-Null
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "10-28-13"

# Row 1: new week's date (Oct 28, 2013) -- keep existing style/format.
$newSheet.Range("B1").Value = 40113

# Row 4: first task for the new week.
$newSheet.Range("A4").Value2 = "Merge GPS and Wifi Code"
$newSheet.Range("B4").Value = 40106
$newSheet.Range("C4").Value = 40107
$newSheet.Range("D4").Value = 1
$newSheet.Range("E4").Value = 3
$newSheet.Range("M4").ClearContents()
$newSheet.Range("N4").ClearContents()

# Row 5: second task for the new week.
$newSheet.Range("A5").Value2 = "Create Code for moving to GPS coordinate"
$newSheet.Range("B5").Value = 40111
$newSheet.Range("C5").ClearContents()
$newSheet.Range("D5").Value = 0.1
$newSheet.Range("E5").Value = 0.5

# The previously-active sheet is no longer the selected tab; move its
# selection off the old "Total time" cell onto A4.
$prevSheet.Activate() | Out-Null
$prevSheet.Range("A4").Select() | Out-Null

# The new sheet becomes the active / selected tab.
$newSheet.Activate() | Out-Null
$newSheet.Range("J4").Select() | Out-Null
